$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{ Row=2; D="96.582.54"; E="  +0.12%  "; Numeric=$false },
  @{ Row=3; D="3.687.00"; E="  -0.11%  "; Numeric=$false },
  @{ Row=4; D="2.44"; E="  +30.86%  "; Numeric=$true },
  @{ Row=5; D="1.00"; E="  +0.03%  "; Numeric=$true },
  @{ Row=6; D="228.48"; E="  -3.22%  "; Numeric=$true },
  @{ Row=7; D="650.84"; E="  +0.19%  "; Numeric=$true },
  @{ Row=8; D="0.437"; E="  +2.40%  "; Numeric=$true },
  @{ Row=9; D="1.15"; E="  +9.46%  "; Numeric=$true },
  @{ Row=10; D="1.00"; E="  +0.00%  "; Numeric=$true },
  @{ Row=11; D="3.685.36"; E="  -0.10%  "; Numeric=$false },
  @{ Row=12; D="47.51"; E="  +7.47%  "; Numeric=$true },
  @{ Row=13; D="0.210"; E="  +2.64%  "; Numeric=$true },
  @{ Row=14; D="0.0000299"; E="  -4.65%  "; Numeric=$true },
  @{ Row=15; D="6.61"; E="  -1.55%  "; Numeric=$true },
  @{ Row=16; D="4.378.05"; E="  -0.01%  "; Numeric=$false },
  @{ Row=17; D="96.374.48"; E="  +0.13%  "; Numeric=$false },
  @{ Row=18; D="8.88"; E="  +1.02%  "; Numeric=$true },
  @{ Row=19; D="3.667.65"; E="  -0.54%  "; Numeric=$false },
  @{ Row=20; D="19.47"; E="  +4.37%  "; Numeric=$true },
  @{ Row=21; D="12.84"; E="  -0.61%  "; Numeric=$true },
  @{ Row=22; D="0.545"; E="  +9.00%  "; Numeric=$true },
  @{ Row=23; D="531.00"; E="  +2.63%  "; Numeric=$true },
  @{ Row=24; D="3.31"; E="  -1.72%  "; Numeric=$true },
  @{ Row=25; D="0.251"; E="  +43.93%  "; Numeric=$true },
  @{ Row=26; D="120.63"; E="  +19.92%  "; Numeric=$true },
  @{ Row=27; D="0.0000209"; E="  +0.27%  "; Numeric=$true },
  @{ Row=28; D="6.82"; E="  -1.46%  "; Numeric=$true },
  @{ Row=29; D="3.890.04"; E="  -0.06%  "; Numeric=$false },
  @{ Row=30; D="12.92"; E="  -1.37%  "; Numeric=$true },
  @{ Row=31; D="13.25"; E="  +9.60%  "; Numeric=$true },
  @{ Row=32; D="2.98"; E="  -0.36%  "; Numeric=$true },
  @{ Row=33; D="1.00"; E="  -0.04%  "; Numeric=$true },
  @{ Row=34; D="0.186"; E="  +0.80%  "; Numeric=$true },
  @{ Row=35; D="33.14"; E="  +3.37%  "; Numeric=$true },
  @{ Row=36; D="1.81"; E="  -2.25%  "; Numeric=$true },
  @{ Row=37; D="0.996"; E="  -0.15%  "; Numeric=$true },
  @{ Row=38; D="0.609"; E="  +3.92%  "; Numeric=$true },
  @{ Row=39; D="605.81"; E="  -7.13%  "; Numeric=$true },
  @{ Row=40; D="1.00"; E="  +0.00%  "; Numeric=$true },
  @{ Row=41; D="8.41"; E="  -4.00%  "; Numeric=$true },
  @{ Row=42; D="7.11"; E="  +3.04%  "; Numeric=$true },
  @{ Row=43; D="0.516"; E="  +21.42%  "; Numeric=$true },
  @{ Row=44; D="0.162"; E="  +2.27%  "; Numeric=$true },
  @{ Row=45; D="0.0501"; E="  +11.78%  "; Numeric=$true },
  @{ Row=46; D="40.34"; E="  +0.10%  "; Numeric=$true },
  @{ Row=47; D="1.99"; E="  -4.13%  "; Numeric=$true },
  @{ Row=48; D="0.966"; E="  +1.45%  "; Numeric=$true },
  @{ Row=49; D="8.97"; E="  +6.28%  "; Numeric=$true },
  @{ Row=50; D="2.28"; E="  +0.87%  "; Numeric=$true },
  @{ Row=51; D="23.53"; E="  -0.14%  "; Numeric=$true }
)

foreach ($item in $updates) {
    $dCell = "D" + $item.Row
    $eCell = "E" + $item.Row
    if ($item.Numeric) {
        $ws.Range($dCell).NumberFormat = "@"
        $ws.Range($dCell).Value = $item.D
        $ws.Range($dCell).Style = "Normal"
    } else {
        $ws.Range($dCell).Value = $item.D
    }
    $ws.Range($eCell).Value = $item.E
}
